$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83:154 down to 84:155.
$ws.Rows.Item(83).Insert()

# Copy the (now shifted-down) row 84 values into the new row 83,
# since the new row at 83 duplicates the prior content of that row
# except for the date in column D which changes.
$src = 84
$dst = 83

$ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($src, 1).Value2
$ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($src, 2).Value2
$ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($src, 3).Value2
$ws.Cells.Item($dst, 4).Value = 44669
$ws.Cells.Item($dst, 4).NumberFormat = $ws.Cells.Item($src, 4).NumberFormat
$ws.Cells.Item($dst, 5).Value = $ws.Cells.Item($src, 5).Value2
$ws.Cells.Item($dst, 6).Value = $ws.Cells.Item($src, 6).Value2
$ws.Cells.Item($dst, 7).Value = $ws.Cells.Item($src, 7).Value2
$ws.Cells.Item($dst, 8).Value = $ws.Cells.Item($src, 8).Value2
$ws.Cells.Item($dst, 9).Value = $ws.Cells.Item($src, 9).Value2
$ws.Cells.Item($dst, 10).Value = $ws.Cells.Item($src, 10).Value2
$ws.Cells.Item($dst, 11).Value = $ws.Cells.Item($src, 11).Value2
$ws.Cells.Item($dst, 12).Value = $ws.Cells.Item($src, 12).Value2
$ws.Cells.Item($dst, 13).Value = $ws.Cells.Item($src, 13).Value2
$ws.Cells.Item($dst, 14).Value = $ws.Cells.Item($src, 14).Value2
$ws.Cells.Item($dst, 15).Value = $ws.Cells.Item($src, 15).Value2
$ws.Cells.Item($dst, 16).Value = $ws.Cells.Item($src, 16).Value2
$ws.Cells.Item($dst, 17).Value = $ws.Cells.Item($src, 17).Value2
$ws.Cells.Item($dst, 18).Value = $ws.Cells.Item($src, 18).Value2
